$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.921.84"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("D3").Value = "1.815.00"
$ws.Range("E3").Value = "  +0.34%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'309.19"
$ws.Range("E5").Value = "  -0.45%  "
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").Value = "'0.4648"
$ws.Range("E7").Value = "  +0.31%  "
$ws.Range("D8").Value = "'0.3657"
$ws.Range("E8").Value = "  -1.61%  "
$ws.Range("D9").Value = "'0.07366"
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("D10").Value = "'0.8683"
$ws.Range("E10").Value = "  -0.68%  "
$ws.Range("D11").Value = "'20.24"
$ws.Range("E11").Value = "  -1.08%  "
$ws.Range("D12").Value = "1.839.46"
$ws.Range("E12").Value = "  -0.67%  "
$ws.Range("D13").Value = "'5.376"
$ws.Range("E13").Value = "  +0.16%  "
$ws.Range("D14").Value = "'0.07107"
$ws.Range("E14").Value = "  +0.90%  "
$ws.Range("D15").Value = "'6.501"
$ws.Range("E15").Value = "  -0.40%  "
$ws.Range("D16").Value = "'91.14"
$ws.Range("E16").Value = "  -1.53%  "
$ws.Range("D17").Value = "'1.003"
$ws.Range("E17").Value = "  +0.13%  "
$ws.Range("E18").Value = "  -0.30%  "
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("D20").Value = "'14.63"
$ws.Range("E20").Value = "  -0.46%  "
$ws.Range("D21").Value = "26.940.20"
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("D22").Value = "'5.289"
$ws.Range("E22").Value = "  -0.42%  "
$ws.Range("E23").Value = "  -0.66%  "
$ws.Range("D24").Value = "2.059.79"
$ws.Range("E24").Value = "  -0.44%  "
$ws.Range("D25").Value = "'1.895"
$ws.Range("E25").Value = "  -0.34%  "
$ws.Range("D26").Value = "'150.73"
$ws.Range("E26").Value = "  -0.66%  "
$ws.Range("D27").Value = "'18.31"
$ws.Range("E27").Value = "  -0.43%  "
$ws.Range("D28").Value = "'2.124"
$ws.Range("E28").Value = "  -0.90%  "
$ws.Range("D29").Value = "'5.250"
$ws.Range("E29").Value = "  -0.61%  "
$ws.Range("D30").Value = "'115.49"
$ws.Range("D31").Value = "'0.08909"
$ws.Range("E31").Value = "  -0.29%  "
$ws.Range("D32").Value = "'0.7531"
$ws.Range("E32").Value = "  -0.38%  "
$ws.Range("D33").Value = "'1.160"
$ws.Range("E33").Value = "  +0.34%  "
$ws.Range("D34").Value = "'4.479"
$ws.Range("E34").Value = "  +0.45%  "
$ws.Range("D35").Value = "'2.902"
$ws.Range("E35").Value = "  -0.90%  "
$ws.Range("E36").Value = "  +0.14%  "
$ws.Range("D37").Value = "'1.095"
$ws.Range("E37").Value = "  -0.91%  "
$ws.Range("D38").Value = "'0.05279"
$ws.Range("E38").Value = "  +0.70%  "
$ws.Range("D39").Value = "'0.01944"
$ws.Range("E39").Value = "  -1.76%  "
$ws.Range("D40").Value = "'2.978"
$ws.Range("E40").Value = "  +1.96%  "
$ws.Range("D41").Value = "'7.223"
$ws.Range("E41").Value = "  +0.34%  "
$ws.Range("E42").Value = "  -0.42%  "
$ws.Range("D43").Value = "'2.301"
$ws.Range("E43").Value = "  -5.25%  "
$ws.Range("D44").Value = "'0.1651"
$ws.Range("E44").Value = "  -0.84%  "
$ws.Range("D45").Value = "'8.402"
$ws.Range("E45").Value = "  -1.32%  "
$ws.Range("D46").Value = "'0.4846"
$ws.Range("E46").Value = "  -2.72%  "
$ws.Range("D47").Value = "'10.44"
$ws.Range("E47").Value = "  +0.45%  "
$ws.Range("E48").Value = "  +0.13%  "
$ws.Range("D49").Value = "'103.15"
$ws.Range("E49").Value = "  -0.80%  "
$ws.Range("D50").Value = "'1.658"
$ws.Range("E50").Value = "  -1.04%  "
$ws.Range("E51").Value = "  -0.14%  "
